$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Suspended Sediment Concentration" (period 5) row; all rows
# below it shift up by one (former row 13 becomes row 12, etc.), and the
# sheet gains a brand-new "QMCI period 15" row at the (new) end, which is
# added further down after the bulk value refresh.
$ws.Rows(12).Delete()

# Refresh every data row (2-24 existing + new 25) with the May-2024 values.
$ws.Range("A2").Value = 'Waikawa at North Manakau Road'
$ws.Range("B2").Value = 'Chlorophyll A'
$ws.Range("C2").Value = 5
$ws.Range("D2").Value = $true
$ws.Range("E2").Value = 'ok'
$ws.Range("F2").Value = 0.952401201978179
$ws.Range("G2").Value = 0.0169491525423729
$ws.Range("H2").Value = 0.796610169491525
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 3.1
$ws.Range("K2").Value = -0.501717032967033
$ws.Range("L2").Value = -1.11695707532747
$ws.Range("M2").Value = -0.0332953826381671
$ws.Range("N2").Value = -16.1844204182914
$ws.Range("O2").Value = 'RepSite'
$ws.Range("P2").Value = 'Extremely likely improving'
$ws.Range("Q2").Value = 1788683
$ws.Range("R2").Value = 5491286
$ws.Range("S2").Value = 'Horowhenua District'
$ws.Range("T2").Value = 'Waiopehu'
$ws.Range("U2").Value = 'Waikawa'
$ws.Range("V2").Value = 'West_9a'
$ws.Range("W2").Value = 'mg/m2'
$ws.Range("A3").Value = 'Waikawa at North Manakau Road'
$ws.Range("B3").Value = 'Visual Clarity'
$ws.Range("C3").Value = 5
$ws.Range("D3").Value = $false
$ws.Range("E3").Value = 'ok'
$ws.Range("F3").Value = 0.999910533883363
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0.8545454545454541
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 3.37
$ws.Range("K3").Value = 0.443517857142857
$ws.Range("L3").Value = 0.277395462797742
$ws.Range("M3").Value = 0.622250249219758
$ws.Range("N3").Value = 13.1607672742688
$ws.Range("O3").Value = 'RepSite'
$ws.Range("P3").Value = 'Virtually certain improving'
$ws.Range("Q3").Value = 1788683
$ws.Range("R3").Value = 5491286
$ws.Range("S3").Value = 'Horowhenua District'
$ws.Range("T3").Value = 'Waiopehu'
$ws.Range("U3").Value = 'Waikawa'
$ws.Range("V3").Value = 'West_9a'
$ws.Range("W3").Value = 'm'
$ws.Range("A4").Value = 'Waikawa at North Manakau Road'
$ws.Range("B4").Value = 'Dissolved Oxygen Concentration'
$ws.Range("C4").Value = 5
$ws.Range("D4").Value = $true
$ws.Range("E4").Value = 'ok'
$ws.Range("F4").Value = 0.469067134141493
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0.9090909090909089
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 10.77
$ws.Range("K4").Value = -0.0062534587123403
$ws.Range("L4").Value = -0.0695027465464326
$ws.Range("M4").Value = 0.0588530407453826
$ws.Range("N4").Value = -0.0580636834943397
$ws.Range("O4").Value = 'RepSite'
$ws.Range("P4").Value = 'As likely as not increasing'
$ws.Range("Q4").Value = 1788683
$ws.Range("R4").Value = 5491286
$ws.Range("S4").Value = 'Horowhenua District'
$ws.Range("T4").Value = 'Waiopehu'
$ws.Range("U4").Value = 'Waikawa'
$ws.Range("V4").Value = 'West_9a'
$ws.Range("W4").Value = 'g/m3'
$ws.Range("A5").Value = 'Waikawa at North Manakau Road'
$ws.Range("B5").Value = 'Dissolved Reactive Phosphorus'
$ws.Range("C5").Value = 5
$ws.Range("D5").Value = $false
$ws.Range("E5").Value = 'WARNING: Sen slope based on tied non-censored values'
$ws.Range("F5").Value = 0.393139708573204
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0.25
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0.011
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = -0.0002731862378459
$ws.Range("M5").Value = 0.0004149773969244
$ws.Range("N5").Value = 0
$ws.Range("O5").Value = 'RepSite'
$ws.Range("P5").Value = 'As likely as not improving'
$ws.Range("Q5").Value = 1788683
$ws.Range("R5").Value = 5491286
$ws.Range("S5").Value = 'Horowhenua District'
$ws.Range("T5").Value = 'Waiopehu'
$ws.Range("U5").Value = 'Waikawa'
$ws.Range("V5").Value = 'West_9a'
$ws.Range("W5").Value = 'mg/L'
$ws.Range("A6").Value = 'Waikawa at North Manakau Road'
$ws.Range("B6").Value = 'E. coli'
$ws.Range("C6").Value = 5
$ws.Range("D6").Value = $true
$ws.Range("E6").Value = 'WARNING: Sen slope influenced by censored values'
$ws.Range("F6").Value = 0.173137135937223
$ws.Range("G6").Value = 0.0357142857142857
$ws.Range("H6").Value = 0.607142857142857
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 46
$ws.Range("K6").Value = 2.00686813186813
$ws.Range("L6").Value = -0.654605715635677
$ws.Range("M6").Value = 5.76795208808083
$ws.Range("N6").Value = 4.36275680840898
$ws.Range("O6").Value = 'RepSite'
$ws.Range("P6").Value = 'Unlikely improving'
$ws.Range("Q6").Value = 1788683
$ws.Range("R6").Value = 5491286
$ws.Range("S6").Value = 'Horowhenua District'
$ws.Range("T6").Value = 'Waiopehu'
$ws.Range("U6").Value = 'Waikawa'
$ws.Range("V6").Value = 'West_9a'
$ws.Range("W6").Value = 'E. coli/100 mL'
$ws.Range("A7").Value = 'Waikawa at North Manakau Road'
$ws.Range("B7").Value = 'Ammoniacal Nitrogen (NH4)'
$ws.Range("C7").Value = 5
$ws.Range("D7").Value = $false
$ws.Range("E7").Value = 'WARNING: Sen slope based on two censored values'
$ws.Range("F7").Value = 0.886095619229102
$ws.Range("G7").Value = 0.903846153846154
$ws.Range("H7").Value = 0.134615384615385
$ws.Range("I7").Value = 2
$ws.Range("J7").Value = 0.005
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 0
$ws.Range("O7").Value = 'RepSite'
$ws.Range("P7").Value = 'Likely improving'
$ws.Range("Q7").Value = 1788683
$ws.Range("R7").Value = 5491286
$ws.Range("S7").Value = 'Horowhenua District'
$ws.Range("T7").Value = 'Waiopehu'
$ws.Range("U7").Value = 'Waikawa'
$ws.Range("V7").Value = 'West_9a'
$ws.Range("W7").Value = 'mg/L'
$ws.Range("A8").Value = 'Waikawa at North Manakau Road'
$ws.Range("B8").Value = 'Nitrite Nitrogen (NO2)'
$ws.Range("C8").Value = 5
$ws.Range("D8").Value = $false
$ws.Range("E8").Value = 'WARNING: Sen slope influenced by censored values'
$ws.Range("F8").Value = 0.998543877095603
$ws.Range("G8").Value = 0.732142857142857
$ws.Range("H8").Value = 0.107142857142857
$ws.Range("I8").Value = 1
$ws.Range("J8").Value = 0.001
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 0
$ws.Range("N8").Value = 0
$ws.Range("O8").Value = 'RepSite'
$ws.Range("P8").Value = 'Virtually certain improving'
$ws.Range("Q8").Value = 1788683
$ws.Range("R8").Value = 5491286
$ws.Range("S8").Value = 'Horowhenua District'
$ws.Range("T8").Value = 'Waiopehu'
$ws.Range("U8").Value = 'Waikawa'
$ws.Range("V8").Value = 'West_9a'
$ws.Range("W8").Value = 'mg/L'
$ws.Range("A9").Value = 'Waikawa at North Manakau Road'
$ws.Range("B9").Value = 'Nitrate Nitrogen (NO3)'
$ws.Range("C9").Value = 5
$ws.Range("D9").Value = $true
$ws.Range("E9").Value = 'ok'
$ws.Range("F9").Value = 0.438865520750858
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 0.857142857142857
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0.06950000000000001
$ws.Range("K9").Value = 0.0005017170329669999
$ws.Range("L9").Value = -0.0035936162727782
$ws.Range("M9").Value = 0.0065641653981205
$ws.Range("N9").Value = 0.721895011463357
$ws.Range("O9").Value = 'RepSite'
$ws.Range("P9").Value = 'As likely as not improving'
$ws.Range("Q9").Value = 1788683
$ws.Range("R9").Value = 5491286
$ws.Range("S9").Value = 'Horowhenua District'
$ws.Range("T9").Value = 'Waiopehu'
$ws.Range("U9").Value = 'Waikawa'
$ws.Range("V9").Value = 'West_9a'
$ws.Range("W9").Value = 'mg/L'
$ws.Range("A10").Value = 'Waikawa at North Manakau Road'
$ws.Range("B10").Value = 'pH'
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = $false
$ws.Range("E10").Value = 'ok'
$ws.Range("F10").Value = 0.033090465300496
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 0.8
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 7.53
$ws.Range("K10").Value = -0.0449261992619925
$ws.Range("L10").Value = -0.08181537083520229
$ws.Range("M10").Value = -0.0067136555920255
$ws.Range("N10").Value = -0.596629472270816
$ws.Range("O10").Value = 'RepSite'
$ws.Range("P10").Value = 'Extremely unlikely increasing'
$ws.Range("Q10").Value = 1788683
$ws.Range("R10").Value = 5491286
$ws.Range("S10").Value = 'Horowhenua District'
$ws.Range("T10").Value = 'Waiopehu'
$ws.Range("U10").Value = 'Waikawa'
$ws.Range("V10").Value = 'West_9a'
$ws.Range("W10").Value = ''
$ws.Range("A11").Value = 'Waikawa at North Manakau Road'
$ws.Range("B11").Value = 'SIN (Soluble Inorganic nitrogen)'
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = $true
$ws.Range("E11").Value = 'ok'
$ws.Range("F11").Value = 0.5
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 0.910714285714286
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0.0785
$ws.Range("K11").Value = 0.0001661737943585
$ws.Range("L11").Value = -0.004526615940258
$ws.Range("M11").Value = 0.0059725282802602
$ws.Range("N11").Value = 0.211686362240137
$ws.Range("O11").Value = 'RepSite'
$ws.Range("P11").Value = 'As likely as not improving'
$ws.Range("Q11").Value = 1788683
$ws.Range("R11").Value = 5491286
$ws.Range("S11").Value = 'Horowhenua District'
$ws.Range("T11").Value = 'Waiopehu'
$ws.Range("U11").Value = 'Waikawa'
$ws.Range("V11").Value = 'West_9a'
$ws.Range("W11").Value = 'g/m3'
$ws.Range("A12").Value = 'Waikawa at North Manakau Road'
$ws.Range("B12").Value = 'Total Nitrogen'
$ws.Range("C12").Value = 5
$ws.Range("D12").Value = $true
$ws.Range("E12").Value = 'WARNING: Sen slope based on tied non-censored values'
$ws.Range("F12").Value = 0.289599497087188
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 0.267857142857143
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0.14
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = -0.0025071600998372
$ws.Range("M12").Value = 0.009955537521962401
$ws.Range("N12").Value = 0
$ws.Range("O12").Value = 'RepSite'
$ws.Range("P12").Value = 'Unlikely improving'
$ws.Range("Q12").Value = 1788683
$ws.Range("R12").Value = 5491286
$ws.Range("S12").Value = 'Horowhenua District'
$ws.Range("T12").Value = 'Waiopehu'
$ws.Range("U12").Value = 'Waikawa'
$ws.Range("V12").Value = 'West_9a'
$ws.Range("W12").Value = 'g/m3'
$ws.Range("A13").Value = 'Waikawa at North Manakau Road'
$ws.Range("B13").Value = 'Total Phosphorus'
$ws.Range("C13").Value = 5
$ws.Range("D13").Value = $false
$ws.Range("E13").Value = 'ok'
$ws.Range("F13").Value = 0.0484600423160867
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0.303571428571429
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0.014
$ws.Range("K13").Value = 0.000431227863046
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 0.001003434065934
$ws.Range("N13").Value = 3.08019902175746
$ws.Range("O13").Value = 'RepSite'
$ws.Range("P13").Value = 'Extremely unlikely improving'
$ws.Range("Q13").Value = 1788683
$ws.Range("R13").Value = 5491286
$ws.Range("S13").Value = 'Horowhenua District'
$ws.Range("T13").Value = 'Waiopehu'
$ws.Range("U13").Value = 'Waikawa'
$ws.Range("V13").Value = 'West_9a'
$ws.Range("W13").Value = 'g/m3'
$ws.Range("A14").Value = 'Waikawa at North Manakau Road'
$ws.Range("B14").Value = 'Turbidity'
$ws.Range("C14").Value = 5
$ws.Range("D14").Value = $false
$ws.Range("E14").Value = 'ok'
$ws.Range("F14").Value = 0.99999716658227
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 0.857142857142857
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0.6850000000000001
$ws.Range("K14").Value = -0.140990989942261
$ws.Range("L14").Value = -0.191969793592545
$ws.Range("M14").Value = -0.09324167154984531
$ws.Range("N14").Value = -20.5826262689432
$ws.Range("O14").Value = 'RepSite'
$ws.Range("P14").Value = 'Virtually certain improving'
$ws.Range("Q14").Value = 1788683
$ws.Range("R14").Value = 5491286
$ws.Range("S14").Value = 'Horowhenua District'
$ws.Range("T14").Value = 'Waiopehu'
$ws.Range("U14").Value = 'Waikawa'
$ws.Range("V14").Value = 'West_9a'
$ws.Range("W14").Value = 'NTU/FNU'
$ws.Range("A15").Value = 'Waikawa at North Manakau Road'
$ws.Range("B15").Value = 'Chlorophyll A'
$ws.Range("C15").Value = 10
$ws.Range("D15").Value = $true
$ws.Range("E15").Value = 'ok'
$ws.Range("F15").Value = 0.342607902638626
$ws.Range("G15").Value = 0.0086206896551724
$ws.Range("H15").Value = 0.698275862068966
$ws.Range("I15").Value = 1
$ws.Range("J15").Value = 3.45
$ws.Range("K15").Value = 0.0550376712328768
$ws.Range("L15").Value = -0.112122552498973
$ws.Range("M15").Value = 0.166694877353359
$ws.Range("N15").Value = 1.59529481834425
$ws.Range("O15").Value = 'RepSite'
$ws.Range("P15").Value = 'As likely as not improving'
$ws.Range("Q15").Value = 1788683
$ws.Range("R15").Value = 5491286
$ws.Range("S15").Value = 'Horowhenua District'
$ws.Range("T15").Value = 'Waiopehu'
$ws.Range("U15").Value = 'Waikawa'
$ws.Range("V15").Value = 'West_9a'
$ws.Range("W15").Value = 'mg/m2'
$ws.Range("A16").Value = 'Waikawa at North Manakau Road'
$ws.Range("B16").Value = 'Chlorophyll A'
$ws.Range("C16").Value = 15
$ws.Range("D16").Value = $true
$ws.Range("E16").Value = 'ok'
$ws.Range("F16").Value = 0.0011059060994492
$ws.Range("G16").Value = 0.0058479532163742
$ws.Range("H16").Value = 0.625730994152047
$ws.Range("I16").Value = 1
$ws.Range("J16").Value = 2.95
$ws.Range("K16").Value = 0.146910540915395
$ws.Range("L16").Value = 0.0552293763347877
$ws.Range("M16").Value = 0.228930058328058
$ws.Range("N16").Value = 4.9800183361151
$ws.Range("O16").Value = 'RepSite'
$ws.Range("P16").Value = 'Exceptionally unlikely improving'
$ws.Range("Q16").Value = 1788683
$ws.Range("R16").Value = 5491286
$ws.Range("S16").Value = 'Horowhenua District'
$ws.Range("T16").Value = 'Waiopehu'
$ws.Range("U16").Value = 'Waikawa'
$ws.Range("V16").Value = 'West_9a'
$ws.Range("W16").Value = 'mg/m2'
$ws.Range("A17").Value = 'Waikawa at North Manakau Road'
$ws.Range("B17").Value = 'ASPM (Macroinvertebrate Average Score Per Metric)'
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = $false
$ws.Range("E17").Value = 'ok'
$ws.Range("F17").Value = 0.768783636774762
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = 1
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0.64
$ws.Range("K17").Value = 0.0106894213366431
$ws.Range("L17").Value = -0.0133745771235129
$ws.Range("M17").Value = 0.0292153749457159
$ws.Range("N17").Value = 1.67022208385048
$ws.Range("O17").Value = 'RepSite'
$ws.Range("P17").Value = 'Likely improving'
$ws.Range("Q17").Value = 1788683
$ws.Range("R17").Value = 5491286
$ws.Range("S17").Value = 'Horowhenua District'
$ws.Range("T17").Value = 'Waiopehu'
$ws.Range("U17").Value = 'Waikawa'
$ws.Range("V17").Value = 'West_9a'
$ws.Range("W17").Value = ''
$ws.Range("A18").Value = 'Waikawa at North Manakau Road'
$ws.Range("B18").Value = 'MCI (Macroinvertebrate Community Index)'
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = $false
$ws.Range("E18").Value = 'ok'
$ws.Range("F18").Value = 0.5
$ws.Range("G18").Value = 0
$ws.Range("H18").Value = 1
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 129.17
$ws.Range("K18").Value = 0.06401512371566211
$ws.Range("L18").Value = -15.4959029944148
$ws.Range("M18").Value = 6.74737150768196
$ws.Range("N18").Value = 0.0495588168426586
$ws.Range("O18").Value = 'RepSite'
$ws.Range("P18").Value = 'As likely as not improving'
$ws.Range("Q18").Value = 1788683
$ws.Range("R18").Value = 5491286
$ws.Range("S18").Value = 'Horowhenua District'
$ws.Range("T18").Value = 'Waiopehu'
$ws.Range("U18").Value = 'Waikawa'
$ws.Range("V18").Value = 'West_9a'
$ws.Range("W18").Value = ''
$ws.Range("A19").Value = 'Waikawa at North Manakau Road'
$ws.Range("B19").Value = 'QMCI (Quantitative Macroinvertebrate Community Index)'
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = $false
$ws.Range("E19").Value = 'ok'
$ws.Range("F19").Value = 0.889664319040077
$ws.Range("G19").Value = 0
$ws.Range("H19").Value = 1
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 7.595
$ws.Range("K19").Value = 0.0500471292527081
$ws.Range("L19").Value = -0.142133761566861
$ws.Range("M19").Value = 0.168681436574581
$ws.Range("N19").Value = 0.658948377257513
$ws.Range("O19").Value = 'RepSite'
$ws.Range("P19").Value = 'Likely improving'
$ws.Range("Q19").Value = 1788683
$ws.Range("R19").Value = 5491286
$ws.Range("S19").Value = 'Horowhenua District'
$ws.Range("T19").Value = 'Waiopehu'
$ws.Range("U19").Value = 'Waikawa'
$ws.Range("V19").Value = 'West_9a'
$ws.Range("W19").Value = ''
$ws.Range("A20").Value = 'Waikawa at North Manakau Road'
$ws.Range("B20").Value = 'ASPM (Macroinvertebrate Average Score Per Metric)'
$ws.Range("C20").Value = 10
$ws.Range("D20").Value = $false
$ws.Range("E20").Value = 'ok'
$ws.Range("F20").Value = 0.141565435331173
$ws.Range("G20").Value = 0
$ws.Range("H20").Value = 1
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0.6645
$ws.Range("K20").Value = -0.0111697247706422
$ws.Range("L20").Value = -0.0188475732030093
$ws.Range("M20").Value = 0.0031728459872031
$ws.Range("N20").Value = -1.68092171115759
$ws.Range("O20").Value = 'RepSite'
$ws.Range("P20").Value = 'Unlikely improving'
$ws.Range("Q20").Value = 1788683
$ws.Range("R20").Value = 5491286
$ws.Range("S20").Value = 'Horowhenua District'
$ws.Range("T20").Value = 'Waiopehu'
$ws.Range("U20").Value = 'Waikawa'
$ws.Range("V20").Value = 'West_9a'
$ws.Range("W20").Value = ''
$ws.Range("A21").Value = 'Waikawa at North Manakau Road'
$ws.Range("B21").Value = 'MCI (Macroinvertebrate Community Index)'
$ws.Range("C21").Value = 10
$ws.Range("D21").Value = $false
$ws.Range("E21").Value = 'ok'
$ws.Range("F21").Value = 0.5
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 1
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 129.585
$ws.Range("K21").Value = -0.0928222596448269
$ws.Range("L21").Value = -1.08352342615014
$ws.Range("M21").Value = 2.1947825633649
$ws.Range("N21").Value = -0.0716304044795516
$ws.Range("O21").Value = 'RepSite'
$ws.Range("P21").Value = 'As likely as not improving'
$ws.Range("Q21").Value = 1788683
$ws.Range("R21").Value = 5491286
$ws.Range("S21").Value = 'Horowhenua District'
$ws.Range("T21").Value = 'Waiopehu'
$ws.Range("U21").Value = 'Waikawa'
$ws.Range("V21").Value = 'West_9a'
$ws.Range("W21").Value = ''
$ws.Range("A22").Value = 'Waikawa at North Manakau Road'
$ws.Range("B22").Value = 'QMCI (Quantitative Macroinvertebrate Community Index)'
$ws.Range("C22").Value = 10
$ws.Range("D22").Value = $false
$ws.Range("E22").Value = 'ok'
$ws.Range("F22").Value = 0.8584345646688269
$ws.Range("G22").Value = 0
$ws.Range("H22").Value = 1
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 7.535
$ws.Range("K22").Value = 0.0224276315789472
$ws.Range("L22").Value = -0.0297498672519984
$ws.Range("M22").Value = 0.07944729919246581
$ws.Range("N22").Value = 0.2976460727133
$ws.Range("O22").Value = 'RepSite'
$ws.Range("P22").Value = 'Likely improving'
$ws.Range("Q22").Value = 1788683
$ws.Range("R22").Value = 5491286
$ws.Range("S22").Value = 'Horowhenua District'
$ws.Range("T22").Value = 'Waiopehu'
$ws.Range("U22").Value = 'Waikawa'
$ws.Range("V22").Value = 'West_9a'
$ws.Range("W22").Value = ''
$ws.Range("A23").Value = 'Waikawa at North Manakau Road'
$ws.Range("B23").Value = 'ASPM (Macroinvertebrate Average Score Per Metric)'
$ws.Range("C23").Value = 15
$ws.Range("D23").Value = $false
$ws.Range("E23").Value = 'ok'
$ws.Range("F23").Value = 0.707964816527247
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 1
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0.6555
$ws.Range("K23").Value = 0.0030918453724604
$ws.Range("L23").Value = -0.0061249835339435
$ws.Range("M23").Value = 0.0120602650987777
$ws.Range("N23").Value = 0.471677402358578
$ws.Range("O23").Value = 'RepSite'
$ws.Range("P23").Value = 'Likely improving'
$ws.Range("Q23").Value = 1788683
$ws.Range("R23").Value = 5491286
$ws.Range("S23").Value = 'Horowhenua District'
$ws.Range("T23").Value = 'Waiopehu'
$ws.Range("U23").Value = 'Waikawa'
$ws.Range("V23").Value = 'West_9a'
$ws.Range("W23").Value = ''
$ws.Range("A24").Value = 'Waikawa at North Manakau Road'
$ws.Range("B24").Value = 'MCI (Macroinvertebrate Community Index)'
$ws.Range("C24").Value = 15
$ws.Range("D24").Value = $false
$ws.Range("E24").Value = 'ok'
$ws.Range("F24").Value = 0.275819582543025
$ws.Range("G24").Value = 0
$ws.Range("H24").Value = 0.866666666666667
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 130
$ws.Range("K24").Value = -0.322944297082228
$ws.Range("L24").Value = -0.830291257475859
$ws.Range("M24").Value = 0.818577091764361
$ws.Range("N24").Value = -0.248418690063252
$ws.Range("O24").Value = 'RepSite'
$ws.Range("P24").Value = 'Unlikely improving'
$ws.Range("Q24").Value = 1788683
$ws.Range("R24").Value = 5491286
$ws.Range("S24").Value = 'Horowhenua District'
$ws.Range("T24").Value = 'Waiopehu'
$ws.Range("U24").Value = 'Waikawa'
$ws.Range("V24").Value = 'West_9a'
$ws.Range("W24").Value = ''
$ws.Range("A25").Value = 'Waikawa at North Manakau Road'
$ws.Range("B25").Value = 'QMCI (Quantitative Macroinvertebrate Community Index)'
$ws.Range("C25").Value = 15
$ws.Range("D25").Value = $false
$ws.Range("E25").Value = 'ok'
$ws.Range("F25").Value = 0.863219294137626
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 1
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 7.535
$ws.Range("K25").Value = 0.0310898652786282
$ws.Range("L25").Value = -0.0182053806112978
$ws.Range("M25").Value = 0.0892456877502847
$ws.Range("N25").Value = 0.412606042184847
$ws.Range("O25").Value = 'RepSite'
$ws.Range("P25").Value = 'Likely improving'
$ws.Range("Q25").Value = 1788683
$ws.Range("R25").Value = 5491286
$ws.Range("S25").Value = 'Horowhenua District'
$ws.Range("T25").Value = 'Waiopehu'
$ws.Range("U25").Value = 'Waikawa'
$ws.Range("V25").Value = 'West_9a'
$ws.Range("W25").Value = ''
